$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "迭代次数为：6，x的值为19.989917407141622"
$ws.Range("A88").Value = "迭代次数为：75，x的值为4.795128926340377"
$ws.Range("A134").Value = "迭代次数为：27，x的值为1.365230011360733"
$ws.Range("A145").Value = "迭代次数为：10，x的值为1.3652300135614255"
$ws.Range("A160").Value = "迭代次数为：3，x的值为-2.5001972959407546"
